$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Insert 3 new rows just above the old closing/thick-bottom row (row 216),
#    so it becomes row 219 again, and copy its neighbour's formatting onto
#    the freshly inserted rows 216-218.
# ---------------------------------------------------------------------------
$ws.Rows.Item(216).Resize(3).Insert()
$ws.Range("B215:C215").Copy()
$ws.Range("B216:C218").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 2) Populate the new "Merk Dagang" / "Agama" / "Unit Kuantitas" UPDATE rows
#    plus the new CREATE rows higher up, writing the values in the same
#    order the original authoring tool produced them in sharedStrings.xml so
#    that duplicate-detection lines up the same way.
# ---------------------------------------------------------------------------

# transaction.update.master.setTradeMark (row 218)
$ws.Range("B218").Value = "transaction.update.master.setTradeMark"
$ws.Range("C218").Value = "Memutakhirkan Data Merk Dagang"

# transaction.create.master.setTradeMark (row 28)
$ws.Range("B28").Value = "transaction.create.master.setTradeMark"
$ws.Range("C28").Value = "Menyimpan Data Baru Merk Dagang"

# transaction.create.master.setReligion (row 27)
$ws.Range("B27").Value = "transaction.create.master.setReligion"
$ws.Range("C27").Value = "Menyimpan Data Baru Agama"

# transaction.update.master.setReligion (row 217)
$ws.Range("B217").Value = "transaction.update.master.setReligion"
$ws.Range("C217").Value = "Memutakhirkan Data Agama"

# transaction.update.master.setQuantityUnit (row 216)
$ws.Range("B216").Value = "transaction.update.master.setQuantityUnit"
$ws.Range("C216").Value = "Memutakhirkan Data Unit Kuantitas"

# transaction.create.master.setQuantityUnit (row 26)
$ws.Range("B26").Value = "transaction.create.master.setQuantityUnit"
$ws.Range("C26").Value = "Menyimpan Data Baru Unit Kuantitas"

# ---------------------------------------------------------------------------
# 3) Update the selection / scroll state to match the post-edit view
#    (pane stays split the same way: freeze at column B / row 3, only the
#    scrolled-to corner and the active cell move).
# ---------------------------------------------------------------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 15
$win.ScrollColumn = 2
$ws.Range("B26").Select()
